$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing X3/Y3 values on the existing row 3
$ws.Range("X3").Value = -0.34999899999999684
$ws.Range("Y3").Value = "Down"

# Copy row 3's formatting down to the two new rows first, so the
# date-format (col A) and percent-format (col S/T) styles carry over
# without Excel minting brand-new number-format styles.
$ws.Range("A3:Y3").Copy($ws.Range("A4:Y4"))
$ws.Range("A3:Y3").Copy($ws.Range("A5:Y5"))

# Row 4 - new trade record
$ws.Range("A4").Value = 42649.61215277778
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 42
$ws.Range("E4").Value = 5125
$ws.Range("F4").Value = 787
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = 31
$ws.Range("I4").Value = 88
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 8771
$ws.Range("L4").Value = 110
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 23
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = "Noun"
$ws.Range("Q4").Value = 35.483823948801813
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.0965
$ws.Range("T4").Value = 0.0269
$ws.Range("U4").Value = 4.82
$ws.Range("V4").Value = 2.2799999999999998
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = -0.34999899999999684
$ws.Range("Y4").Value = "Down"

# Row 5 - new trade record
$ws.Range("A5").Value = 42649.63553240741
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 52
$ws.Range("E5").Value = 5930
$ws.Range("F5").Value = 1047
$ws.Range("G5").Value = 67
$ws.Range("H5").Value = 31
$ws.Range("I5").Value = 88
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 10520
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 65
$ws.Range("O5").Value = 8
$ws.Range("P5").Value = "Noun"
$ws.Range("Q5").Value = 35.483823948801813
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.0965
$ws.Range("T5").Value = 0.0269
$ws.Range("U5").Value = 4.82
$ws.Range("V5").Value = 2.2799999999999998
$ws.Range("W5").Value = 0

# Row 5 has no X/Y entries in the source data - clear what the format
# copy brought over so those cells stay blank.
$ws.Range("X5:Y5").ClearContents()

